$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy cell formatting from existing template cells ---
$ws.Range("A1").Copy()
$ws.Range("A40:B40").PasteSpecial(-4122)
$ws.Range("A59:B59").PasteSpecial(-4122)
$ws.Range("A77:B77").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("A42:B42").PasteSpecial(-4122)
$ws.Range("A61:B61").PasteSpecial(-4122)
$ws.Range("A79:B79").PasteSpecial(-4122)
$ws.Range("A4").Copy()
$ws.Range("A43:B43").PasteSpecial(-4122)
$ws.Range("A62:B62").PasteSpecial(-4122)
$ws.Range("A80:B80").PasteSpecial(-4122)
$ws.Range("A5").Copy()
$ws.Range("A44").PasteSpecial(-4122)
$ws.Range("A63").PasteSpecial(-4122)
$ws.Range("A81").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B44").PasteSpecial(-4122)
$ws.Range("B63").PasteSpecial(-4122)
$ws.Range("B81").PasteSpecial(-4122)
$ws.Range("A6").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$ws.Range("A46").PasteSpecial(-4122)
$ws.Range("A47").PasteSpecial(-4122)
$ws.Range("A48").PasteSpecial(-4122)
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("A65").PasteSpecial(-4122)
$ws.Range("A66").PasteSpecial(-4122)
$ws.Range("A67").PasteSpecial(-4122)
$ws.Range("A82").PasteSpecial(-4122)
$ws.Range("A83").PasteSpecial(-4122)
$ws.Range("A84").PasteSpecial(-4122)
$ws.Range("A85").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("B45").PasteSpecial(-4122)
$ws.Range("B46").PasteSpecial(-4122)
$ws.Range("B47").PasteSpecial(-4122)
$ws.Range("B48").PasteSpecial(-4122)
$ws.Range("F54:G54").PasteSpecial(-4122)
$ws.Range("I54:L54").PasteSpecial(-4122)
$ws.Range("F55:G55").PasteSpecial(-4122)
$ws.Range("I55:L55").PasteSpecial(-4122)
$ws.Range("F56:G56").PasteSpecial(-4122)
$ws.Range("I56:L56").PasteSpecial(-4122)
$ws.Range("B64").PasteSpecial(-4122)
$ws.Range("B65").PasteSpecial(-4122)
$ws.Range("B66").PasteSpecial(-4122)
$ws.Range("B67").PasteSpecial(-4122)
$ws.Range("F72:G72").PasteSpecial(-4122)
$ws.Range("I72:L72").PasteSpecial(-4122)
$ws.Range("F73:G73").PasteSpecial(-4122)
$ws.Range("I73:L73").PasteSpecial(-4122)
$ws.Range("F74:G74").PasteSpecial(-4122)
$ws.Range("I74:L74").PasteSpecial(-4122)
$ws.Range("B82").PasteSpecial(-4122)
$ws.Range("B83").PasteSpecial(-4122)
$ws.Range("B84").PasteSpecial(-4122)
$ws.Range("B85").PasteSpecial(-4122)
$ws.Range("F90:G90").PasteSpecial(-4122)
$ws.Range("I90:L90").PasteSpecial(-4122)
$ws.Range("F91:G91").PasteSpecial(-4122)
$ws.Range("I91:L91").PasteSpecial(-4122)
$ws.Range("F92:G92").PasteSpecial(-4122)
$ws.Range("I92:L92").PasteSpecial(-4122)
$ws.Range("A12").Copy()
$ws.Range("A51:M51").PasteSpecial(-4122)
$ws.Range("A69:M69").PasteSpecial(-4122)
$ws.Range("A87:M87").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("A53:L53").PasteSpecial(-4122)
$ws.Range("A71:L71").PasteSpecial(-4122)
$ws.Range("A89:L89").PasteSpecial(-4122)
$ws.Range("A15").Copy()
$ws.Range("A54").PasteSpecial(-4122)
$ws.Range("A55").PasteSpecial(-4122)
$ws.Range("A56").PasteSpecial(-4122)
$ws.Range("A72").PasteSpecial(-4122)
$ws.Range("A73").PasteSpecial(-4122)
$ws.Range("A74").PasteSpecial(-4122)
$ws.Range("A90").PasteSpecial(-4122)
$ws.Range("A91").PasteSpecial(-4122)
$ws.Range("A92").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("B54:C54").PasteSpecial(-4122)
$ws.Range("H54").PasteSpecial(-4122)
$ws.Range("B55:C55").PasteSpecial(-4122)
$ws.Range("H55").PasteSpecial(-4122)
$ws.Range("B56:C56").PasteSpecial(-4122)
$ws.Range("H56").PasteSpecial(-4122)
$ws.Range("B72:C72").PasteSpecial(-4122)
$ws.Range("H72").PasteSpecial(-4122)
$ws.Range("B73:C73").PasteSpecial(-4122)
$ws.Range("H73").PasteSpecial(-4122)
$ws.Range("B74:C74").PasteSpecial(-4122)
$ws.Range("H74").PasteSpecial(-4122)
$ws.Range("B90:C90").PasteSpecial(-4122)
$ws.Range("H90").PasteSpecial(-4122)
$ws.Range("B91:C91").PasteSpecial(-4122)
$ws.Range("H91").PasteSpecial(-4122)
$ws.Range("B92:C92").PasteSpecial(-4122)
$ws.Range("H92").PasteSpecial(-4122)
$ws.Range("D15").Copy()
$ws.Range("D54:E54").PasteSpecial(-4122)
$ws.Range("D55:E55").PasteSpecial(-4122)
$ws.Range("D56:E56").PasteSpecial(-4122)
$ws.Range("D72:E72").PasteSpecial(-4122)
$ws.Range("D73:E73").PasteSpecial(-4122)
$ws.Range("D74:E74").PasteSpecial(-4122)
$ws.Range("D90:E90").PasteSpecial(-4122)
$ws.Range("D91:E91").PasteSpecial(-4122)
$ws.Range("D92:E92").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Set cell values ---
$ws.Range("A40").Value = 'API Summary'
$ws.Range("A42").Value = 'API NAME'
$ws.Range("A43").Value = 'IssueIssueByMonth'
$ws.Range("A44").Value = 'Description'
$ws.Range("B44").Value = 'The function of this api is to collect total issue in perticular Month in perticular project'
$ws.Range("A45").Value = 'Total Test Case Executed'
$ws.Range("B45").Value = 3
$ws.Range("A46").Value = 'Pending'
$ws.Range("B46").Value = 0
$ws.Range("A47").Value = 'Failed'
$ws.Range("B47").Value = 0
$ws.Range("A48").Value = 'Passed'
$ws.Range("B48").Value = 3
$ws.Range("A51").Value = 'Test case Summary'
$ws.Range("A53").Value = 'Test case No'
$ws.Range("B53").Value = 'Name'
$ws.Range("C53").Value = 'Function '
$ws.Range("D53").Value = '% TC Executed'
$ws.Range("E53").Value = '% TC Pending'
$ws.Range("F53").Value = 'Priority'
$ws.Range("G53").Value = 'Remark'
$ws.Range("H53").Value = 'Input to test'
$ws.Range("I53").Value = 'Response Code Expected'
$ws.Range("J53").Value = 'Response Code Got'
$ws.Range("K53").Value = 'Response Content Expected'
$ws.Range("L53").Value = 'Response Content Expected'
$ws.Range("A54").Value = 'TCIBM01'
$ws.Range("B54").Value = 'get_number_of_issue_by_month_
api_testing_by_giving_right_parameter'
$ws.Range("C54").Value = 'This Test case is use to validate 
data inserted by user is right '
$ws.Range("D54").Value = 1
$ws.Range("E54").Value = 0
$ws.Range("F54").Value = 'High'
$ws.Range("G54").Value = 'No Error'
$ws.Range("H54").Value = '{"project_id":202}'
$ws.Range("I54").Value = 200
$ws.Range("J54").Value = 200
$ws.Range("K54").Value = '{''Issue'': 4}'
$ws.Range("L54").Value = '{''Issue'': 4}'
$ws.Range("A55").Value = 'TCIBM02'
$ws.Range("B55").Value = 'get_number_of_issue_by_month_api
_testing_by_giving_missing_parameter'
$ws.Range("C55").Value = 'This test case is designed to validate 
the wheather user is given numerical 
values instead of string'
$ws.Range("D55").Value = 1
$ws.Range("E55").Value = 0
$ws.Range("F55").Value = 'High'
$ws.Range("G55").Value = 'No Error'
$ws.Range("H55").Value = '{}'
$ws.Range("I55").Value = 400
$ws.Range("J55").Value = 400
$ws.Range("K55").Value = '{''error'': ''bad values''}'
$ws.Range("L55").Value = '{''error'': ''bad values''}'
$ws.Range("A56").Value = 'TCIBM03'
$ws.Range("B56").Value = 'get_number_of_issue_by_month_api
_testing_by_giving_Wrong_datatype_parameter'
$ws.Range("C56").Value = 'This test case is designed 
to validate the wheather user is given 
numerical values instead of string'
$ws.Range("D56").Value = 1
$ws.Range("E56").Value = 0
$ws.Range("F56").Value = 'Hight'
$ws.Range("G56").Value = 'No Error'
$ws.Range("H56").Value = '{"project_id":"*&^%$"}'
$ws.Range("I56").Value = 400
$ws.Range("J56").Value = 400
$ws.Range("K56").Value = '{''Error'': ''Wrong data type of project id''}'
$ws.Range("L56").Value = '{''Error'': ''Wrong data type of project id''}'
$ws.Range("A59").Value = 'API Summary'
$ws.Range("A61").Value = 'API NAME'
$ws.Range("A62").Value = 'IssueIssueByWeek'
$ws.Range("A63").Value = 'Description'
$ws.Range("B63").Value = 'The function of this api is to collect total issue in perticular Week in perticular project'
$ws.Range("A64").Value = 'Total Test Case Executed'
$ws.Range("B64").Value = 3
$ws.Range("A65").Value = 'Pending'
$ws.Range("B65").Value = 0
$ws.Range("A66").Value = 'Failed'
$ws.Range("B66").Value = 0
$ws.Range("A67").Value = 'Passed'
$ws.Range("B67").Value = 3
$ws.Range("A69").Value = 'Test case Summary'
$ws.Range("A71").Value = 'Test case No'
$ws.Range("B71").Value = 'Name'
$ws.Range("C71").Value = 'Function '
$ws.Range("D71").Value = '% TC Executed'
$ws.Range("E71").Value = '% TC Pending'
$ws.Range("F71").Value = 'Priority'
$ws.Range("G71").Value = 'Remark'
$ws.Range("H71").Value = 'Input to test'
$ws.Range("I71").Value = 'Response Code Expected'
$ws.Range("J71").Value = 'Response Code Got'
$ws.Range("K71").Value = 'Response Content Expected'
$ws.Range("L71").Value = 'Response Content Expected'
$ws.Range("A72").Value = 'TCIBW01'
$ws.Range("B72").Value = 'get_number_of_issue_by_week_
api_testing_by_giving_right_parameter'
$ws.Range("C72").Value = 'This Test case is use to validate 
data inserted by user is right '
$ws.Range("D72").Value = 1
$ws.Range("E72").Value = 0
$ws.Range("F72").Value = 'High'
$ws.Range("G72").Value = 'No Error'
$ws.Range("H72").Value = '{"project_id":202}'
$ws.Range("I72").Value = 200
$ws.Range("J72").Value = 200
$ws.Range("K72").Value = '{''Issue'': 4}'
$ws.Range("L72").Value = '{''Issue'': 4}'
$ws.Range("A73").Value = 'TCIBW02'
$ws.Range("B73").Value = 'get_number_of_issue_by_week_api
_testing_by_giving_missing_parameter'
$ws.Range("C73").Value = 'This test case is designed to validate 
the wheather user is given numerical 
values instead of string'
$ws.Range("D73").Value = 1
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 'High'
$ws.Range("G73").Value = 'No Error'
$ws.Range("H73").Value = '{}'
$ws.Range("I73").Value = 400
$ws.Range("J73").Value = 400
$ws.Range("K73").Value = '{''error'': ''bad values''}'
$ws.Range("L73").Value = '{''error'': ''bad values''}'
$ws.Range("A74").Value = 'TCIBW03'
$ws.Range("B74").Value = 'get_number_of_issue_by_week_api
_testing_by_giving_Wrong_datatype_parameter'
$ws.Range("C74").Value = 'This test case is designed 
to validate the wheather user is given 
numerical values instead of string'
$ws.Range("D74").Value = 1
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 'Hight'
$ws.Range("G74").Value = 'No Error'
$ws.Range("H74").Value = '{"project_id":"*&^%$"}'
$ws.Range("I74").Value = 400
$ws.Range("J74").Value = 400
$ws.Range("K74").Value = '{''Error'': ''Wrong data type of project id''}'
$ws.Range("L74").Value = '{''Error'': ''Wrong data type of project id''}'
$ws.Range("A77").Value = 'API Summary'
$ws.Range("A79").Value = 'API NAME'
$ws.Range("A80").Value = 'IssueIssueByQuarter'
$ws.Range("A81").Value = 'Description'
$ws.Range("B81").Value = 'The function of this api is to collect total issue in perticular quarter in perticular project'
$ws.Range("A82").Value = 'Total Test Case Executed'
$ws.Range("B82").Value = 3
$ws.Range("A83").Value = 'Pending'
$ws.Range("B83").Value = 0
$ws.Range("A84").Value = 'Failed'
$ws.Range("B84").Value = 0
$ws.Range("A85").Value = 'Passed'
$ws.Range("B85").Value = 3
$ws.Range("A87").Value = 'Test case Summary'
$ws.Range("A89").Value = 'Test case No'
$ws.Range("B89").Value = 'Name'
$ws.Range("C89").Value = 'Function '
$ws.Range("D89").Value = '% TC Executed'
$ws.Range("E89").Value = '% TC Pending'
$ws.Range("F89").Value = 'Priority'
$ws.Range("G89").Value = 'Remark'
$ws.Range("H89").Value = 'Input to test'
$ws.Range("I89").Value = 'Response Code Expected'
$ws.Range("J89").Value = 'Response Code Got'
$ws.Range("K89").Value = 'Response Content Expected'
$ws.Range("L89").Value = 'Response Content Expected'
$ws.Range("A90").Value = 'TCIBQ01'
$ws.Range("B90").Value = 'get_number_of_issue_by_Quarter_
api_testing_by_giving_right_parameter'
$ws.Range("C90").Value = 'This Test case is use to validate 
data inserted by user is right '
$ws.Range("D90").Value = 1
$ws.Range("E90").Value = 0
$ws.Range("F90").Value = 'High'
$ws.Range("G90").Value = 'No Error'
$ws.Range("H90").Value = '{"project_id":202}'
$ws.Range("I90").Value = 200
$ws.Range("J90").Value = 200
$ws.Range("K90").Value = '{''Issue'': 4}'
$ws.Range("L90").Value = '{''Issue'': 4}'
$ws.Range("A91").Value = 'TCIBQ02'
$ws.Range("B91").Value = 'get_number_of_issue_by_Quarter_api
_testing_by_giving_missing_parameter'
$ws.Range("C91").Value = 'This test case is designed to validate 
the wheather user is given numerical 
values instead of string'
$ws.Range("D91").Value = 1
$ws.Range("E91").Value = 0
$ws.Range("F91").Value = 'High'
$ws.Range("G91").Value = 'No Error'
$ws.Range("H91").Value = '{}'
$ws.Range("I91").Value = 400
$ws.Range("J91").Value = 400
$ws.Range("K91").Value = '{''error'': ''bad values''}'
$ws.Range("L91").Value = '{''error'': ''bad values''}'
$ws.Range("A92").Value = 'TCIBQ03'
$ws.Range("B92").Value = 'get_number_of_issue_by_Quarter_api
_testing_by_giving_Wrong_datatype_parameter'
$ws.Range("C92").Value = 'This test case is designed 
to validate the wheather user is given 
numerical values instead of string'
$ws.Range("D92").Value = 1
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 'Hight'
$ws.Range("G92").Value = 'No Error'
$ws.Range("H92").Value = '{"project_id":"*&^%$"}'
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 400
$ws.Range("K92").Value = '{''Error'': ''Wrong data type of project id''}'
$ws.Range("L92").Value = '{''Error'': ''Wrong data type of project id''}'

# --- Row heights (approximate natural wrap height) ---
$ws.Rows.Item(44).RowHeight = 45.0
$ws.Rows.Item(54).RowHeight = 30.0
$ws.Rows.Item(55).RowHeight = 45.0
$ws.Rows.Item(56).RowHeight = 45.0
$ws.Rows.Item(63).RowHeight = 30.0
$ws.Rows.Item(72).RowHeight = 30.0
$ws.Rows.Item(73).RowHeight = 45.0
$ws.Rows.Item(74).RowHeight = 45.0
$ws.Rows.Item(81).RowHeight = 30.0
$ws.Rows.Item(90).RowHeight = 30.0
$ws.Rows.Item(91).RowHeight = 45.0
$ws.Rows.Item(92).RowHeight = 45.0

# --- Merged cells ---
$ws.Range("A40:B40").Merge()
$ws.Range("A42:B42").Merge()
$ws.Range("A43:B43").Merge()
$ws.Range("A51:M51").Merge()
$ws.Range("A59:B59").Merge()
$ws.Range("A61:B61").Merge()
$ws.Range("A62:B62").Merge()
$ws.Range("A69:M69").Merge()
$ws.Range("A77:B77").Merge()
$ws.Range("A79:B79").Merge()
$ws.Range("A80:B80").Merge()
$ws.Range("A87:M87").Merge()

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 51.166666666666664
$ws.Columns.Item(11).ColumnWidth = 35.3
$ws.Columns.Item(12).ColumnWidth = 34.8

# --- View / selection state ---
$ws.Range("D80").Select()

Write-Output "done"